$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date text in B1
$ws.Range("B1").Value = "18/03/2023"

# Update B/C numeric values for rows 2-20
$data = @{
    2  = @(338, 38)
    3  = @(185, 38)
    4  = @(4, 38)
    5  = @(116, 38)
    6  = @(20, 38)
    7  = @(40, 38)
    8  = @(21, 8)
    9  = @(48, 18)
    10 = @(198, 38)
    11 = @(127, 38)
    12 = @(256.1, 38)
    13 = @(263, 38)
    14 = @(344, 38)
    15 = @(102, 38)
    16 = @(92, 28)
    17 = @(57, 38)
    18 = @(2, 38)
    19 = @(5, 38)
    20 = @(22, 38)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
}
